$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 282.5164233333333
$ws.Range("H2").Value = 847.54927
$ws.Range("I2").Value = 0.7504954445259187
$ws.Range("J2").Value = 0.7504954445259185
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 96.88352066666666
$ws.Range("N2").Value = 290.650562
$ws.Range("O2").Value = 0.4395029568526832
$ws.Range("P2").Value = 0.4395029568526832
$ws.Range("Q2").Value = 27371.18573868774
$ws.Range("R2").Value = 246340.6716481897
$ws.Range("S2").Value = 0.3298449669736101
$ws.Range("T2").Value = 0.3298449669736101

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 282.5164233333333
$ws.Range("H3").Value = 847.54927
$ws.Range("I3").Value = 0.7504954445259187
$ws.Range("J3").Value = 0.7504954445259185
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 80.72275766666667
$ws.Range("N3").Value = 242.168273
$ws.Range("O3").Value = 0.3661911792188718
$ws.Range("P3").Value = 0.3661911792188718
$ws.Range("Q3").Value = 22805.50477759008
$ws.Range("R3").Value = 205249.5429983107
$ws.Range("S3").Value = 0.2748248118293375
$ws.Range("T3").Value = 0.2748248118293375

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 282.5164233333333
$ws.Range("H4").Value = 847.54927
$ws.Range("I4").Value = 0.7504954445259187
$ws.Range("J4").Value = 0.7504954445259185
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 42.83255866666667
$ws.Range("N4").Value = 128.497676
$ws.Range("O4").Value = 0.194305863928445
$ws.Range("P4").Value = 0.194305863928445
$ws.Range("Q4").Value = 12100.90127672184
$ws.Range("R4").Value = 108908.1114904965
$ws.Range("S4").Value = 0.145825665722971
$ws.Range("T4").Value = 0.145825665722971

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 68.514867
$ws.Range("H5").Value = 205.544601
$ws.Range("I5").Value = 0.1820074562714184
$ws.Range("J5").Value = 0.1820074562714184
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 96.88352066666666
$ws.Range("N5").Value = 290.650562
$ws.Range("O5").Value = 0.4395029568526832
$ws.Range("P5").Value = 0.4395029568526832
$ws.Range("Q5").Value = 6637.961532968417
$ws.Range("R5").Value = 59741.65379671576
$ws.Range("S5").Value = 0.07999281520052383
$ws.Range("T5").Value = 0.07999281520052381

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 68.514867
$ws.Range("H6").Value = 205.544601
$ws.Range("I6").Value = 0.1820074562714184
$ws.Range("J6").Value = 0.1820074562714184
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 80.72275766666667
$ws.Range("N6").Value = 242.168273
$ws.Range("O6").Value = 0.3661911792188718
$ws.Range("P6").Value = 0.3661911792188718
$ws.Range("Q6").Value = 5530.709005404897
$ws.Range("R6").Value = 49776.38104864407
$ws.Range("S6").Value = 0.06664952503865795
$ws.Range("T6").Value = 0.06664952503865794

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 68.514867
$ws.Range("H7").Value = 205.544601
$ws.Range("I7").Value = 0.1820074562714184
$ws.Range("J7").Value = 0.1820074562714184
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 42.83255866666667
$ws.Range("N7").Value = 128.497676
$ws.Range("O7").Value = 0.194305863928445
$ws.Range("P7").Value = 0.194305863928445
$ws.Range("Q7").Value = 2934.667060316364
$ws.Range("R7").Value = 26412.00354284728
$ws.Range("S7").Value = 0.03536511603223663
$ws.Range("T7").Value = 0.03536511603223663

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 25.40860066666666
$ws.Range("H8").Value = 76.22580199999999
$ws.Range("I8").Value = 0.06749709920266306
$ws.Range("J8").Value = 0.06749709920266304
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 96.88352066666666
$ws.Range("N8").Value = 290.650562
$ws.Range("O8").Value = 0.4395029568526832
$ws.Range("P8").Value = 0.4395029568526832
$ws.Range("Q8").Value = 2461.674687800079
$ws.Range("R8").Value = 22155.07219020072
$ws.Range("S8").Value = 0.0296651746785493
$ws.Range("T8").Value = 0.02966517467854929

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 25.40860066666666
$ws.Range("H9").Value = 76.22580199999999
$ws.Range("I9").Value = 0.06749709920266306
$ws.Range("J9").Value = 0.06749709920266304
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 80.72275766666667
$ws.Range("N9").Value = 242.168273
$ws.Range("O9").Value = 0.3661911792188718
$ws.Range("P9").Value = 0.3661911792188718
$ws.Range("Q9").Value = 2051.052314264438
$ws.Range("R9").Value = 18459.47082837994
$ws.Range("S9").Value = 0.02471684235087636
$ws.Range("T9").Value = 0.02471684235087635

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.40860066666666
$ws.Range("H10").Value = 76.22580199999999
$ws.Range("I10").Value = 0.06749709920266306
$ws.Range("J10").Value = 0.06749709920266304
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 42.83255866666667
$ws.Range("N10").Value = 128.497676
$ws.Range("O10").Value = 0.194305863928445
$ws.Range("P10").Value = 0.194305863928445
$ws.Range("Q10").Value = 1088.315378692906
$ws.Range("R10").Value = 9794.838408236152
$ws.Range("S10").Value = 0.0131150821732374
$ws.Range("T10").Value = 0.0131150821732374
